$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "Save" column, matching the style of the existing
# header row (column G / "sum")
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save column values per row
$saveValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
